$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 4 and row 5 for columns D, J, K, L, M, N, O, P, Q
# (all other columns are identical between the two rows already).

# --- Save row 4 original values ---
$D4 = $ws.Range("D4").Value2
$J4 = $ws.Range("J4").Value2
$K4 = $ws.Range("K4").Value2
$L4 = $ws.Range("L4").Value2
$M4 = $ws.Range("M4").Value2
$N4 = $ws.Range("N4").Value2
$O4 = $ws.Range("O4").Value2
$P4 = $ws.Range("P4").Value2
$Q4 = $ws.Range("Q4").Value2

# --- Save row 5 original values ---
$D5 = $ws.Range("D5").Value2
$J5 = $ws.Range("J5").Value2
$K5 = $ws.Range("K5").Value2
$L5 = $ws.Range("L5").Value2
$M5 = $ws.Range("M5").Value2
$N5 = $ws.Range("N5").Value2
$O5 = $ws.Range("O5").Value2
$P5 = $ws.Range("P5").Value2
$Q5 = $ws.Range("Q5").Value2

# --- Write row 5's original values into row 4 ---
$ws.Range("D4").Value = $D5
$ws.Range("J4").Value = $J5
$ws.Range("K4").Value = $K5
$ws.Range("L4").Value = $L5
$ws.Range("M4").Value = $M5
$ws.Range("N4").Value = $N5
$ws.Range("O4").Value = $O5
$ws.Range("P4").Value = $P5
$ws.Range("Q4").Value = $Q5

# --- Write row 4's original values into row 5 ---
$ws.Range("D5").Value = $D4
$ws.Range("J5").Value = $J4
$ws.Range("K5").Value = $K4
$ws.Range("L5").Value = $L4
$ws.Range("M5").Value = $M4
$ws.Range("N5").Value = $N4
$ws.Range("O5").Value = $O4
$ws.Range("P5").Value = $P4
$ws.Range("Q5").Value = $Q4
